$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.093421219042667758
$ws.Range("B1").Value = 0.093301433484285212
$ws.Range("A2").Value = -0.047190585555027198
$ws.Range("B2").Value = 0.046854563605322141
$ws.Range("A3").Value = 0.091945117157489875
$ws.Range("B3").Value = -0.092268719589640114
$ws.Range("A4").Value = -0.19172556619417236
$ws.Range("B4").Value = 0.19074064146294845
$ws.Range("A5").Value = -0.18474064177012206
$ws.Range("B5").Value = 0.18275838276295353
$ws.Range("A6").Value = -0.092789645854551406
$ws.Range("B6").Value = 0.092664063850642187
$ws.Range("A7").Value = -0.072664064230878012
$ws.Range("B7").Value = 0.07237461735154227
$ws.Range("A8").Value = -0.052374617735455153
$ws.Range("B8").Value = 0.052159086374526709
$ws.Range("A9").Value = -0.046159086702480145
$ws.Range("B9").Value = 0.045987045013310635
$ws.Range("A10").Value = -0.039987045344894057
$ws.Range("B10").Value = 0.039963986343224178
$ws.Range("A11").Value = -0.035463986668354153
$ws.Range("B11").Value = 0.035429495362915731
$ws.Range("A12").Value = -0.029429495695401986
$ws.Range("B12").Value = 0.029333948108865915
$ws.Range("A13").Value = -0.023333948444664188
$ws.Range("B13").Value = 0.023312652105305531
$ws.Range("A14").Value = -0.027085378964391893
$ws.Range("B14").Value = 0.027053102679827745
$ws.Range("A15").Value = -0.021053103017737662
$ws.Range("B15").Value = 0.021027691658136227
$ws.Range("A16").Value = -0.015027691997196779
$ws.Range("B16").Value = 0.015004508286680984
$ws.Range("A17").Value = -0.009004508627218577
$ws.Range("B17").Value = 0.0089999996462131193
$ws.Range("A18").Value = -0.036110848248664951
$ws.Range("B18").Value = 0.03609674471286084
$ws.Range("A19").Value = -0.027096745028222013
$ws.Range("B19").Value = 0.027013822270560262
$ws.Range("A20").Value = -0.018013822588724082
$ws.Range("B20").Value = 0.018004296675531606
$ws.Range("A21").Value = -0.0090042969941137585
$ws.Range("B21").Value = 0.008999999681112314
$ws.Range("A22").Value = -0.093933796854093643
$ws.Range("B22").Value = 0.093625193723491762
$ws.Range("A23").Value = -0.084625194043337792
$ws.Range("B23").Value = 0.084125009307058463
$ws.Range("A24").Value = -0.042125009773894107
$ws.Range("B24").Value = 0.041999999530694332
$ws.Range("A25").Value = -0.083968737227220203
$ws.Range("B25").Value = 0.083874181415108495
$ws.Range("A26").Value = -0.077874181734873815
$ws.Range("B26").Value = 0.077759216553591415
$ws.Range("A27").Value = -0.07175921687483644
$ws.Range("B27").Value = 0.071390327815415766
$ws.Range("A28").Value = -0.06539032814242951
$ws.Range("B28").Value = 0.065156603823536052
$ws.Range("A29").Value = -0.062520096160655569
$ws.Range("B29").Value = 0.062164977525215548
$ws.Range("A30").Value = -0.042164977918557778
$ws.Range("B30").Value = 0.042019471292542576
$ws.Range("A31").Value = -0.027019471667570016
$ws.Range("B31").Value = 0.027000840892164035
$ws.Range("A32").Value = -0.0060008412935435018
$ws.Range("B32").Value = 0.0059999996629569452
